$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the duplicated "Author"/"Subject" columns (G:H) that were
#    accidentally left in the sheet.
# ---------------------------------------------------------------------------
$ws.Range("G1:H4").EntireColumn.Delete()

# ---------------------------------------------------------------------------
# 2. Rewrite the table body: rows 3 and 4 swap places (Erica Jong's record
#    now comes before Emily's) and Erica Jong's ExpectedResult is filled in
#    with "success" instead of being left blank.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Author"
$ws.Range("B1").Value = "Subject"
$ws.Range("C1").Value = "Edition"
$ws.Range("D1").Value = "Format"
$ws.Range("E1").Value = "AgeGroup"
$ws.Range("F1").Value = "ExpectedResult"

$ws.Range("A2").Value = "Antony"
$ws.Range("B2").Value = "fiction"
$ws.Range("C2").Value = "Edition 4"
$ws.Range("D2").Value = "NewsPaper"
$ws.Range("E2").Value = "teen"
$ws.Range("F2").Value = "success"

$ws.Range("A3").Value = "Erica Jong"
$ws.Range("B3").Value = "non-fiction"
$ws.Range("C3").Value = "Edition 2"
$ws.Range("D3").Value = "Magazines"
$ws.Range("E3").Value = "adult"
$ws.Range("F3").Value = "success"

$ws.Range("A4").Value = "Emily"
$ws.Range("B4").Value = "horror"
$ws.Range("C4").Value = "Edition 2"
$ws.Range("D4").Value = "Magazines"
$ws.Range("E4").Value = "kids"
$ws.Range("F4").Value = "failure"

# ---------------------------------------------------------------------------
# 3. Formatting touch-ups: the header row no longer uses the taller 28.8pt
#    row height, column widths for E/F were adjusted and the selection
#    moved off-sheet to D11.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).EntireRow.AutoFit()
$ws.Rows.Item(2).EntireRow.AutoFit()
$ws.Rows.Item(3).EntireRow.AutoFit()
$ws.Rows.Item(4).EntireRow.AutoFit()

$ws.Columns.Item(5).ColumnWidth = 12.35
$ws.Columns.Item(6).ColumnWidth = 14.0

$rng = $ws.Range("A1:F4")
$rng.WrapText = $true
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4160

$ws.Range("F3").WrapText = $false
$ws.Range("F3").HorizontalAlignment = -4131
$ws.Range("F3").VerticalAlignment = -4160

$ws.Range("D11").Select()
